$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("방어팀", "방어팀 펫", "방어팀 스순", "공격팀", "공격팀 펫", "공격팀 스순", "속공", "날짜", "상대 길드", "기준"),
  @("프레이야, 오공, 실베스타", "연지", "프2오2프1", "트루드, 겔리두스, 스파이크", "루", "겔2트2스2", "후", 260105, "모현", "수비"),
  @("프레이야, 바네사, 키리엘", "연지", "바1프2키2", "프레이야, 바네사, 유신", "연지", "바1프2유1", "선", 260105, "모현", "수비"),
  @("카일, 카구라, 파이", "이린", "구2카2카1", "트루드, 파이, 엘리시아", "카람", "엘1파1트2", "후", 260105, "모현", "수비"),
  @("프레이야, 바네사, 키리엘", "연지", "바1프2키2", "플라튼, 콜트, 챈슬러", "파이크", "챈2콜1콜2", "후", 260105, "모현", "수비"),
  @("트루드, 오공, 엘리시아", "파이크", "오2엘1트2", "플라튼, 엘리스, 리나", "루", "플2엘2리1", "후", 260105, "모현", "수비"),
  @("트루드, 스파이크, 아멜리아", "파이크", "아1스2트2", "프레이야, 연희, 키리엘", "델로", "프2프1연2", "선", 260105, "모현", "수비"),
  @("카일, 카구라, 파이", "이린", "구2카2카1", "트루드, 겔리두스, 엘리시아", "크리", "겔2엘1트2", "후", 260105, "모현", "수비"),
  @("오공, 겔리두스, 엘리시아", "루", "오2엘2겔2", "오공, 리나, 엘리스", "파이크", "오2엘1엘2", "후", 260105, "모현", "수비"),
  @("오공, 프레이야, 엘리시아", "유", "프2엘1오2", "프레이야, 바네사, 키리엘", "연지", "바1프2프1", "선", 260105, "모현", "수비"),
  @("오공, 겔리두스, 엘리시아", "맬패로", "오2엘2겔2", "오공", "파이크", "오2오1", "후", 260105, "모현", "수비"),
  @("트루드, 겔리두스, 엘리시아", "리첼", "겔2엘2트2", "플라튼, 카론, 루시", "파이크", "플2루2카1", "후", 260105, "모현", "수비"),
  @("카일, 카구라, 콜트", "이린", "구2카1카2", "트루드, 오공, 엘리시아", "루", "트2오2엘1", "후", 260105, "모현", "수비"),
  @("카일, 카구라, 파이", "이린", "구2파1카1", "카일, 카구라, 파이", "이린", "구2카1카2", "선", 260105, "모현", "수비"),
  @("밀리아, 연희, 키리엘", "연지", "키2연2밀2", "오공, 겔리두스, 엘리시아", "파이크", "오2엘1겔1", "후", 260105, "모현", "수비"),
  @("오공, 겔리두스, 에이스", "루", "겔2오2겔1", "풍연, 카일, 파이", "이린", "파1카2파2", "후", 260105, "모현", "수비"),
  @("플라튼, 트루드, 콜트", "파이크", "플2콜1콜2", "오공, 레이첼, 아멜리아", "루", "레2아2오2", "후", 260105, "모현", "수비"),
  @("프레이야, 바네사, 키리엘", "연지", "바1프2키2", "오공, 루디, 리나", "파이크", "오2오1루2", "후", 260105, "모현", "수비"),
  @("카일, 카구라, 파이", "이린", "구2카1카2", "카일, 카구라, 파이", "이린", "파1카1카2", "선", 260105, "모현", "수비"),
  @("오공, 스파이크, 아멜리아", "파이크", "아2오2스2", "오공, 스파이크, 아멜리아", "파이크", "아1아2오2", "선", 260105, "모현", "수비"),
  @("프레이야, 바네사, 키리엘", "유", "바1프2키2", "프레이야, 연희, 바네사", "연지", "바1프2프1", "선", 260105, "모현", "수비"),
  @("프레이야, 오공, 엘리시아", "유", "엘1오2프2", "프레이야, 바네사, 키리엘", "연지", "바1프2키1", "선", 260105, "모현", "수비"),
  @("카일, 카구라, 파이", "이린", "구2카1카2", "카일, 카구라, 파이", "이린", "구2카2카1", "선", 260105, "모현", "수비"),
  @("카일, 카구라, 파이", "이린", "구2카2카1", "트루드, 제이브, 카구라", "카람", "구2제1트2", "후", 260105, "모현", "수비"),
  @("오공, 겔리두스, 엘리시아", "루", "오2겔2오1", "오공, 엘리스, 리나", "파이크", "오2", "후", 260105, "모현", "수비"),
  @("플라튼, 프레이야, 콜트", "루", "프2콜1콜2", "카일, 파이, 아멜리아", "이린", "아2카1카2", "선", 260105, "모현", "수비"),
  @("트루드, 스파이크, 아멜리아", "크리", "아1스2트2", "챈슬러, 스파이크, 엘리스", "파이크", "챈2엘1스2", "후", 260105, "모현", "수비"),
  @("오공, 겔리두스, 에이스", "유", "오2겔2에1", "카일, 파이, 아멜리아", "이린", "아2카1카2", "선", 260105, "모현", "수비"),
  @("밀리아, 프레이야, 바네사", "연지", "바1밀2프2", "카일, 파이, 카구라", "이린", "구2카1카2", "선", 260105, "모현", "수비"),
  @("오공, 스파이크, 아멜리아", "파이크", "오2아1아2", "프레이야, 연희, 키리엘", "연지", "프2연2프1", "선", 260105, "모현", "수비"),
  @("트루드, 스파이크, 아멜리아", "파이크", "아1스2트2", "오공, 스파이크, 겔리두스", "파이크", "오2겔2오1", "선", 260105, "모현", "수비"),
  @("트루드, 엘리시아, 플라튼", "루", "트2엘1플2", "프레이야, 트루드, 플라튼", "연지", "프2트2프1", "선", 260105, "모현", "수비"),
  @("카일, 카구라, 파이", "이린", "카2카1파2", "트루드, 제이브, 카구라", "리첼", "트2구2제1", "선", 260105, "모현", "수비"),
  @("밀리아, 바네사, 키리엘", "맬패로", "밀2바1키2", "오공, 겔리두스, 엘리시아", "파이크", "오2겔2오1", "후", 260105, "모현", "수비"),
  @("트루드, 스파이크, 아멜리아", "루", "아1트2스2", "카일, 카구라, 파이", "이린", "구2카1카2", "선", 260105, "모현", "수비"),
  @("카일, 카구라, 파이", "이린", "파1카1카2", "플라튼, 트루드, 콜트", "루", "트2콜1콜2", "후", 260105, "모현", "수비"),
  @("플라튼, 실베스타, 아멜리아", "루", "예약없음", "카일, 카구라, 아멜리아", "이린", "아2카1카2", "선", 260105, "모현", "수비"),
  @("플라튼, 트루드, 아멜리아", "파이크", "아1아2트2", "프레이야, 키리엘, 오를리", "연지", "프2프1키2", "선", 260105, "모현", "수비"),
  @("트루드, 스파이크, 아멜리아", "루", "아2트2스1", "카일, 카구라, 파이", "이린", "카1파2카2", "선", 260105, "모현", "수비"),
  @("플라튼, 프레이야, 콜트", "루", "프2콜1콜2", "크리스, 로지, 아킬라", "맬패로", "아2크1크2", "후", 260105, "모현", "수비"),
  @("오공, 트루드, 엘리시아", "루", "오2트2엘2", "프레이야, 바네사, 키리엘", "연지", "바1프2키2", "선", 260105, "모현", "수비"),
  @("카일, 카구라, 파이", "이린", "구2카2카1", "트루드, 제이브, 카구라", "이린", "구2제1트2", "선", 260105, "모현", "수비"),
  @("트루드, 스파이크, 아멜리아", "루", "아1트2아2", "카론, 스파이크, 챈슬러", "파이크", "챈2스2카1", "후", 260105, "모현", "수비"),
  @("프레이야, 키리엘, 바네사", "연지", "바1프2프1", "오공, 겔리두스, 스파이크", "루", "겔2오2겔1", "후", 260105, "모현", "수비"),
  @("트루드, 스파이크, 아멜리아", "파이크", "아1아2트2", "프레이야, 바네사, 키리엘", "연지", "프2프1바1", "선", 260105, "모현", "수비"),
  @("카일, 카구라, 파이", "크리", "구2카1카2", "플라튼, 콜트, 쥬리", "루", "쥬2콜1콜2", "후", 260105, "모현", "수비"),
  @("카일, 카구라, 파이", "이린", "카1카2파2", "트루드, 제이브, 카구라", "카람", "제1트2제2", "후", 260105, "모현", "수비"),
  @("밀리아, 프레이야, 키리엘", "연지", "키2프2프1", "오공, 겔리두스, 스파이크", "파이크", "오2겔2스2", "후", 260105, "모현", "수비"),
  @("오공, 겔리두스, 엘리시아", "맬패로", "오2엘2겔2", "카일, 카구라, 아멜리아", "이린", "구2카1카2", "선", 260105, "모현", "수비"),
  @("카일, 카구라, 파이", "크리", "카1구2카2", "트루드, 제이브, 카구라", "헬레핀", "제1트2제2", "후", 260105, "모현", "수비"),
  @("트루드, 엘리시아, 아멜리아", "리첼", "아1아2트2", "프레이야, 연희, 키리엘", "연지", "프2프1연2", "선", 260105, "모현", "수비"),
  @("프레이야, 연희, 바네사", "연지", "바1프2연2", "오공, 루디, 엘리스", "파이크", "루2오2", "후", 260105, "모현", "수비"),
  @("겔리두스, 트루드, 엘리시아", "파이크", "겔1트2엘2", "리, 엘리스, 플라튼", "카람", "리1플2엘1", "후", 260105, "모현", "수비"),
  @("트루드, 엘리시아, 플라튼", "루", "트2엘2플2", "크리스, 녹스, 로지", "맬패로", "크1녹2크2", "후", 260105, "모현", "수비"),
  @("오공, 겔리두스, 스파이크", "파이크", "겔2오2겔1", "오공, 겔리두스, 에이스", "루", "오2에2에1", "후", 260105, "모현", "수비"),
  @("카일, 카구라, 파이", "이린", "구2카2카1", "오공, 겔리두스, 스파이크", "파이크", "오2겔1겔2", "후", 260105, "모현", "수비"),
  @("밀리아, 프레이야, 바네사", "맬패로", "바1프2프1", "오공, 겔리두스, 스파이크", "크리", "오2겔2스2", "후", 260105, "모현", "수비"),
  @("카일, 카구라, 파이", "델로", "구2카1카2", "플라튼, 콜트, 트루드", "루", "콜2트2", "후", 260105, "모현", "수비"),
  @("프레이야, 바네사, 키리엘", "연지", "프2바1프1", "트루드, 콜트, 플라튼", "루", "콜1콜2트2", "선", 260105, "모현", "수비"),
  @("프레이야, 바네사, 키리엘", "연지", "프2바1프1", "프레이야, 바네사, 키리엘", "연지", "바1프2프1", "선", 260105, "모현", "수비"),
  @("밀리아, 바네사, 키리엘", "맬패로", "바1밀2키2", "프레이야, 바네사, 키리엘", "연지", "바1프2프1", "후", 260105, "모현", "수비"),
  @("오공, 엘리시아, 겔리두스", "루", "겔1오2엘1", "트루드, 오공, 스파이크", "파이크", "오2오1트2", "후", 260105, "모현", "수비"),
  @("카일, 카구라, 파이", "이린", "구2카1카2", "카일, 카구라, 파이", "이린", "카1카2파2", "후", 260105, "모현", "수비"),
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 1
    $rowdata = $data[$i]
    for ($j = 0; $j -lt $rowdata.Length; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $rowdata[$j]
    }
}

$ws.Range("A1:J63").HorizontalAlignment = -4108
$ws.Range("A1:J63").VerticalAlignment = -4108

$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("D67").Select()

Write-Output "edit applied"
